$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B,C to C,D)
$ws.Columns("B").Insert()

# Set the new header and query values
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.ethnicity IN ['HISPANIC_OR_LATINO'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match column widths: new column B should match column A's width, column C/D keep prior widths
$ws.Columns("A:B").Width = 75.81640625
$ws.Columns("C").Width = 70.26953125
$ws.Columns("D").Width = 28.54296875

# Apply wrap text style to B2 (matching A2's style)
$ws.Range("B2").WrapText = $true

# Update selection to match diff
$ws.Range("B2").Select()
